$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-66 down to 51-67.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly Mango price record.
$ws.Cells.Item(50, 1).Value  = 11
$ws.Cells.Item(50, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value  = "Bíobío"
$ws.Cells.Item(50, 4).Value  = 44460
$ws.Cells.Item(50, 5).Value  = 8
$ws.Cells.Item(50, 6).Value  = "Fruta"
$ws.Cells.Item(50, 7).Value  = 100108
$ws.Cells.Item(50, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(50, 9).Value  = 100108002
$ws.Cells.Item(50, 10).Value = "Mango"
$ws.Cells.Item(50, 11).Value = "Sin especificar"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 200
$ws.Cells.Item(50, 14).Value = 8000
$ws.Cells.Item(50, 15).Value = 8500
$ws.Cells.Item(50, 16).Value = 8250
$ws.Cells.Item(50, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(50, 18).Value = "Brasil"
$ws.Cells.Item(50, 19).Value = 2062
$ws.Cells.Item(50, 20).Value = 4
